$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price (column D) cells whose new values are plain numeric
# strings need to stay stored as text (matching the sheet convention used for
# every other price cell). Force a Text number format before writing the value
# so Excel does not auto-convert them to the Number type, then restore the
# cell style so no visual/style change is introduced.
$textPriceRows = @(5, 6, 10, 11, 12, 16, 17, 19, 21, 22, 23, 25, 27, 28, 29, 30, 31, 32, 34, 35, 38, 40, 41, 44, 46, 48, 50, 51)
foreach ($r in $textPriceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "46.685.97"
$ws.Range("E2").Value = "  +6.27%  "
$ws.Range("D3").Value = "2.309.74"
$ws.Range("E3").Value = "  +5.31%  "
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").Value = "301.03"
$ws.Range("E5").Value = "  +2.22%  "
$ws.Range("D6").Value = "102.12"
$ws.Range("E6").Value = "  +14.53%  "
$ws.Range("E7").Value = "  +1.72%  "
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("E9").Value = "  +9.52%  "
$ws.Range("D10").Value = "36.75"
$ws.Range("E10").Value = "  +14.61%  "
$ws.Range("D11").Value = "0.0804"
$ws.Range("E11").Value = "  +4.09%  "
$ws.Range("D12").Value = "7.38"
$ws.Range("E12").Value = "  +8.44%  "
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "2.662.60"
$ws.Range("E14").Value = "  +5.41%  "
$ws.Range("D15").Value = "2.306.46"
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("D16").Value = "14.05"
$ws.Range("E16").Value = "  +6.82%  "
$ws.Range("D17").Value = "0.822"
$ws.Range("E17").Value = "  +6.30%  "
$ws.Range("D18").Value = "46.657.21"
$ws.Range("E18").Value = "  +7.01%  "
$ws.Range("D19").Value = "13.20"
$ws.Range("E19").Value = "  +22.00%  "
$ws.Range("D20").Value = "0.0₃0947"
$ws.Range("E20").Value = "  +6.41%  "
$ws.Range("D21").Value = "6.14"
$ws.Range("E21").Value = "  +5.29%  "
$ws.Range("D22").Value = "66.89"
$ws.Range("E22").Value = "  +5.92%  "
$ws.Range("D23").Value = "248.38"
$ws.Range("E23").Value = "  +7.16%  "
$ws.Range("E24").Value = "  +7.39%  "
$ws.Range("D25").Value = "1.96"
$ws.Range("E25").Value = "  +7.27%  "
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").Value = "43.76"
$ws.Range("E27").Value = "  +20.41%  "
$ws.Range("D28").Value = "2.25"
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("D29").Value = "9.96"
$ws.Range("E29").Value = "  +7.67%  "
$ws.Range("D30").Value = "20.16"
$ws.Range("E30").Value = "  +4.87%  "
$ws.Range("D31").Value = "5.80"
$ws.Range("E31").Value = "  +10.02%  "
$ws.Range("D32").Value = "0.0805"
$ws.Range("E32").Value = "  +9.41%  "
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("D34").Value = "2.62"
$ws.Range("E34").Value = "  +4.38%  "
$ws.Range("D35").Value = "3.15"
$ws.Range("E35").Value = "  +11.71%  "
$ws.Range("E36").Value = "  +9.50%  "
$ws.Range("E37").Value = "  +3.86%  "
$ws.Range("D38").Value = "1.83"
$ws.Range("E38").Value = "  +11.25%  "
$ws.Range("E39").Value = "  +20.26%  "
$ws.Range("D40").Value = "4.09"
$ws.Range("E40").Value = "  +16.76%  "
$ws.Range("D41").Value = "3.49"
$ws.Range("E41").Value = "  +13.82%  "
$ws.Range("E42").Value = "  +7.71%  "
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").Value = "1.98"
$ws.Range("E44").Value = "  +18.98%  "
$ws.Range("D45").Value = "1.851.84"
$ws.Range("E45").Value = "  +3.04%  "
$ws.Range("D46").Value = "88.59"
$ws.Range("E46").Value = "  +22.48%  "
$ws.Range("E47").Value = "  +12.64%  "
$ws.Range("D48").Value = "74.80"
$ws.Range("E48").Value = "  +15.42%  "
$ws.Range("E49").Value = "  +11.33%  "
$ws.Range("D50").Value = "97.59"
$ws.Range("E50").Value = "  +6.28%  "
$ws.Range("D51").Value = "54.77"
$ws.Range("E51").Value = "  +10.83%  "

# Restore the original (default) cell style on the cells we touched above so
# the Text number format does not linger as a visible style change.
foreach ($r in $textPriceRows) {
    $ws.Range("D$r").Style = "Normal"
}
